# Add two new columns, I (I0) and J (IF), to Sheet1.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Headers
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Copy the header style/format from H1 (bold, centered, bordered) onto I1:J1
# so they match the rest of the header row.
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)

# Data for rows 2-36: column I is 1 for every row except row 23 (which is 4),
# column J mirrors column H except row 23 (which is 6 instead of 3).
$iVals = @{
    2=1; 3=1; 4=1; 5=1; 6=1; 7=1; 8=1; 9=1; 10=1;
    11=1; 12=1; 13=1; 14=1; 15=1; 16=1; 17=1; 18=1; 19=1; 20=1;
    21=1; 22=1; 23=4; 24=1; 25=1; 26=1; 27=1; 28=1; 29=1; 30=1;
    31=1; 32=1; 33=1; 34=1; 35=1; 36=1
}
$jVals = @{
    2=4; 3=6; 4=5; 5=7; 6=6; 7=5; 8=5; 9=5; 10=5;
    11=6; 12=5; 13=5; 14=5; 15=5; 16=6; 17=7; 18=5; 19=6; 20=5;
    21=6; 22=6; 23=6; 24=4; 25=4; 26=6; 27=7; 28=7; 29=6; 30=5;
    31=4; 32=4; 33=5; 34=5; 35=3; 36=2
}

for ($r = 2; $r -le 36; $r++) {
    $ws.Cells.Item($r, 9).Value = $iVals[$r]
    $ws.Cells.Item($r, 10).Value = $jVals[$r]
}
